$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 35 (this shifts old rows 35-56 down to 37-58)
$ws.Rows.Item(35).Resize(2).Insert()

# Row 35 (new)
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35, 3).Value = "Ñuble"
$ws.Cells.Item(35, 4).Value = 44763
$ws.Cells.Item(35, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(35, 6).Value = 100112040
$ws.Cells.Item(35, 7).Value = "Cilantro"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 200
$ws.Cells.Item(35, 11).Value = 700
$ws.Cells.Item(35, 12).Value = 800
$ws.Cells.Item(35, 13).Value = 750
$ws.Cells.Item(35, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(35, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(35, 16).Value = 750
$ws.Cells.Item(35, 17).Value = 1
$ws.Cells.Item(35, 18).Value = "Hortaliza"

# Row 36 (new)
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 44763
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = 100112040
$ws.Cells.Item(36, 7).Value = "Cilantro"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Segunda"
$ws.Cells.Item(36, 10).Value = 150
$ws.Cells.Item(36, 11).Value = 600
$ws.Cells.Item(36, 12).Value = 600
$ws.Cells.Item(36, 13).Value = 600
$ws.Cells.Item(36, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(36, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(36, 16).Value = 600
$ws.Cells.Item(36, 17).Value = 1
$ws.Cells.Item(36, 18).Value = "Hortaliza"
